$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Clone the formatting of the last existing data row (row 13) onto the new
# row 14 so the new row picks up the same per-column cell styles, then
# overwrite the cell values that differ for the new monster card.
$src = $ws.Range("A13:J13")
$dst = $ws.Range("A14:J14")
$src.Copy($dst)

$ws.Range("A14").Value = 57000011
$ws.Range("F14").Value = "lp2mp"
$ws.Range("B14").Value = "远古龙"

# Grow the table to include the new row.
$tbl.Resize($ws.Range("A3:J14"))

# Keep the active-pane selection in sync with the new last row, matching
# Excel's behaviour after extending the table by one row.
$ws.Range("G14").Select()
